$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "الذهب"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "معلومة"

$ws.Range("G7").Select()
